# Print stmts for visibility and control for blank (np.nan) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Email values (column B, rows 2-6)
$emails = @(
    "aaimzbrgoa@aclBa.mhm",
    "agelnlcesn@dmmor.iam",
    "glicom@mmane.agi",
    "nmmeidw@aogig.lco",
    "yremcomsaiam@gtrij.aal"
)

# New Phone values (column C, rows 2-6)
$phones = @(
    "04-203-5005902",
    "90-005-3245020",
    "40-900-2203550",
    "40-020-5029350",
    "20-050-2095304"
)

# New Address values (column D, rows 2-6)
$addresses = @(
    "yssa/rndmtdr/peeet",
    "pyta/tmersdd/rnsee",
    "smnp/ardettd/eresy",
    "tedp/aystres/rdemn",
    "yend/esedrsm/rtatp"
)

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    Write-Host "Updating row $row"
    $ws.Cells.Item($row, 2).Value = $emails[$i]
    $ws.Cells.Item($row, 3).Value = $phones[$i]
    $ws.Cells.Item($row, 4).Value = $addresses[$i]
}

Write-Host "Done updating Email, Phone, and Address columns."
